# Daily update at 8 AM UTC
# Adds the next day's row (row 31) to the "Wins Over Time" tracking sheet
# and restores the previous last row (row 30) to the standard date/time
# number format, since row 31 is now the final ("today") row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30 was previously the last row and used the short "YYYY-MM-DD" date
# format reserved for the newest entry. Now that a new row is being
# appended, row 30 reverts to the regular "YYYY-MM-DD HH:MM:SS" format.
$ws.Cells.Item(30, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data as row 31.
$ws.Cells.Item(31, 1).Value = 45615
$ws.Cells.Item(31, 2).Value = 77
$ws.Cells.Item(31, 3).Value = 64
$ws.Cells.Item(31, 4).Value = 76

# Row 31 is now the newest entry, so it gets the short date-only format.
$ws.Cells.Item(31, 1).NumberFormat = "YYYY-MM-DD"
